$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.126707792282104
$ws.Range("B1").Value = 0.8723677396774292
$ws.Range("C1").Value = 4.115319728851318
$ws.Range("D1").Value = 2.784535884857178
$ws.Range("E1").Value = 0.784727931022644
